$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = 3.65
$ws.Range("S3").Value = 3.65
$ws.Range("U3").Value = 2.04
$ws.Range("I4").Value = 7.2
$ws.Range("J4").Value = 2.82
$ws.Range("N4").Value = 2.26
$ws.Range("W4").Value = 2
$ws.Range("G6").Value = 2.64
$ws.Range("H6").Value = 3.4
$ws.Range("J6").Value = 2.96
$ws.Range("V6").Value = 1.33
$ws.Range("W6").Value = 1.61
$ws.Range("X6").Value = 970
$ws.Range("Y6").Value = 970
$ws.Range("AB6").Value = 970
$ws.Range("AC6").Value = 970
$ws.Range("AD6").Value = 970
$ws.Range("AF6").Value = 970
$ws.Range("AG6").Value = 970
$ws.Range("K8").Value = 4.7
$ws.Range("U8").Value = 2.38
$ws.Range("V8").Value = 2.1
$ws.Range("I11").Value = 2.2
$ws.Range("J11").Value = 3.35
$ws.Range("G13").Value = 1.96
$ws.Range("I13").Value = 4.7
$ws.Range("M13").Value = 1.03
$ws.Range("N13").Value = 5
$ws.Range("O13").Value = 1.2
$ws.Range("R13").Value = 1.57
$ws.Range("S13").Value = 2.48
$ws.Range("T13").Value = 1.61
$ws.Range("U13").Value = 2.32
$ws.Range("V13").Value = 1.27
$ws.Range("W13").Value = 2.04
$ws.Range("X13").Value = 24
$ws.Range("Y13").Value = 22
$ws.Range("Z13").Value = 36
$ws.Range("AA13").Value = 85
$ws.Range("AB13").Value = 13
$ws.Range("AC13").Value = 10.5
$ws.Range("AD13").Value = 18.5
$ws.Range("AE13").Value = 48
$ws.Range("AF13").Value = 14.5
$ws.Range("AG13").Value = 11.5
$ws.Range("AH13").Value = 17.5
$ws.Range("AI13").Value = 48
$ws.Range("AJ13").Value = 22
$ws.Range("AK13").Value = 18.5
$ws.Range("AL13").Value = 29
$ws.Range("AM13").Value = 70
$ws.Range("AN13").Value = 9.4
$ws.Range("AO13").Value = 42
$ws.Range("K14").Value = 4
$ws.Range("S14").Value = 3.1
$ws.Range("X14").Value = 970
$ws.Range("P15").Value = 1.69
$ws.Range("S15").Value = 3.8
$ws.Range("Z15").Value = 60
$ws.Range("F16").Value = 2.16
$ws.Range("G16").Value = 2.52
$ws.Range("H16").Value = 3.25
$ws.Range("I16").Value = 4.1
$ws.Range("J16").Value = 3.25
$ws.Range("K16").Value = 4.4
$ws.Range("L16").Value = 1.24
$ws.Range("N16").Value = 3.75
$ws.Range("O16").Value = 1.25
$ws.Range("P16").Value = 2.14
$ws.Range("Q16").Value = 1.63
$ws.Range("R16").Value = 1.45
$ws.Range("S16").Value = 2.58
$ws.Range("T16").Value = 1.65
$ws.Range("U16").Value = 2.3
$ws.Range("V16").Value = 1.37
$ws.Range("W16").Value = 1.69
$ws.Range("X16").Value = 22
$ws.Range("Y16").Value = 970
$ws.Range("Z16").Value = 30
$ws.Range("AA16").Value = 70
$ws.Range("AD16").Value = 970
$ws.Range("AE16").Value = 44
$ws.Range("AG16").Value = 13
$ws.Range("AH16").Value = 970
$ws.Range("AI16").Value = 50
$ws.Range("AK16").Value = 25
$ws.Range("AL16").Value = 40
$ws.Range("AM16").Value = 90
$ws.Range("AO16").Value = 36
$ws.Range("J17").Value = 3.85
$ws.Range("L17").Value = 1.28
$ws.Range("M17").Value = 1.06
$ws.Range("N17").Value = 3.9
$ws.Range("O17").Value = 1.28
$ws.Range("R17").Value = 1.39
$ws.Range("S17").Value = 3.1
$ws.Range("T17").Value = 1.75
$ws.Range("U17").Value = 2.12
$ws.Range("X17").Value = 16
$ws.Range("Y17").Value = 17
$ws.Range("Z17").Value = 34
$ws.Range("AB17").Value = 9.800000000000001
$ws.Range("AC17").Value = 8.800000000000001
$ws.Range("AD17").Value = 18
$ws.Range("AF17").Value = 12.5
$ws.Range("AG17").Value = 10.5
$ws.Range("AJ17").Value = 22
$ws.Range("AK17").Value = 20
$ws.Range("AL17").Value = 36
$ws.Range("AN17").Value = 13
$ws.Range("AO17").Value = 60
$ws.Range("H18").Value = 2.26
$ws.Range("L18").Value = 1.45
$ws.Range("N18").Value = 3.65
$ws.Range("V19").Value = 1.73
$ws.Range("F20").Value = 1.88
$ws.Range("I20").Value = 5.5
$ws.Range("K20").Value = 6
$ws.Range("N20").Value = 1.03
$ws.Range("O20").Value = 1.39
$ws.Range("P20").Value = 1.25
$ws.Range("Q20").Value = 2.18
$ws.Range("S20").Value = 2.18
$ws.Range("F21").Value = 4.8
$ws.Range("I21").Value = 1.91
$ws.Range("J21").Value = 3.5
$ws.Range("K21").Value = 3.9
$ws.Range("P21").Value = 1.82
$ws.Range("T21").Value = 1.88
$ws.Range("V21").Value = 2.1
$ws.Range("AA21").Value = 20
$ws.Range("AF21").Value = 40
$ws.Range("AL21").Value = 80
$ws.Range("F22").Value = 1.9
$ws.Range("I22").Value = 5.8
$ws.Range("K22").Value = 3.45
$ws.Range("L22").Value = 1.01
$ws.Range("M22").Value = 1.13
$ws.Range("N22").Value = 2.42
$ws.Range("O22").Value = 1.62
$ws.Range("R22").Value = 1.15
$ws.Range("S22").Value = 6
$ws.Range("T22").Value = 2.34
$ws.Range("U22").Value = 1.51
$ws.Range("V22").Value = 1.2
$ws.Range("W22").Value = 2
$ws.Range("X22").Value = 9.6
$ws.Range("Y22").Value = 13.5
$ws.Range("Z22").Value = 44
$ws.Range("AA22").Value = 230
$ws.Range("AB22").Value = 6
$ws.Range("AC22").Value = 8.4
$ws.Range("AD22").Value = 26
$ws.Range("AE22").Value = 150
$ws.Range("AF22").Value = 10.5
$ws.Range("AG22").Value = 12
$ws.Range("AH22").Value = 36
$ws.Range("AI22").Value = 180
$ws.Range("AJ22").Value = 25
$ws.Range("AK22").Value = 36
$ws.Range("AL22").Value = 75
$ws.Range("AM22").Value = 360
$ws.Range("AN22").Value = 29
$ws.Range("AO22").Value = 290
